# Applies:
#  1. Inserts a new "Player Info" worksheet as the first sheet, containing
#     ID / NAME / BATTING_HAND / BOWL_STYLE columns for player 4926
#     (James Alexander McCollum).
#  2. Renames the "MATCH_CARD_LINK" column in the existing "ODI Batting"
#     sheet to "MATCH_CODE" and replaces the full scorecard URLs with the
#     bare numeric match code that was embedded in them.

$wb = $excel.ActiveWorkbook

# --- 1. New "Player Info" sheet, placed before "ODI Batting" -----------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"
$playerInfo.Move($wb.Worksheets.Item(1))

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Header formatting to match the bold/bordered/centred style used by the
# "ODI Batting" sheet's header row.
$hdr = $playerInfo.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Player 4926 - James Alexander McCollum. ID is stored as text (leading
# apostrophe), matching the original sheet's convention of keeping
# numeric-looking identifiers as text.
$playerInfo.Range("A2").Value = "'4926"
$playerInfo.Range("B2").Value = "James Alexander McCollum"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium"

# --- 2. Update "ODI Batting" sheet --------------------------------------
$odi = $wb.Worksheets.Item("ODI Batting")
$odi.Range("D1").Value = "MATCH_CODE"

$matchCodes = @{
    2  = "4257"
    3  = "4291"
    4  = "4295"
    5  = "4299"
    6  = "4301"
    7  = "4343"
    8  = "4347"
    9  = "4352"
    10 = "4397"
    11 = "4448"
}

foreach ($row in $matchCodes.Keys) {
    $odi.Cells.Item($row, 4).Value = "'" + $matchCodes[$row]
}
